# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) switch from the deck's custom
#    "Table_0" style to the built-in "No Style, No Grid" table style.
# 2) The presentation's theme colour palette switches from the
#    "Red Violet" / Integral palette to the standard "Office" palette
#    (this is what the slide master / presentation-wide theme part
#    ends up containing after the edit).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1. Re-style the three tables.
# ---------------------------------------------------------------
$newTableStyleId = "{C14F183F-6B72-4C44-B8DB-13A6CD03B206}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------
# 2. Swap the theme colour scheme to the "Office" palette.
# ---------------------------------------------------------------
function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($idx = 1; $idx -le $officeColors.Length; $idx++) {
    $themeColors.Colors($idx).RGB = HexToVbaRgb $officeColors[$idx - 1]
}
